$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numCues for heavensComingDown (row 3): 39 -> 40
$ws.Range("B3").Value = 40

# Update songNotes for heavensComingDown (row 3): "quad heavens gate" -> "2 x double heavens gate"
$ws.Range("C3").Value = "~cue 9: 2 x double heavens gate                                                                                            ~cues 10, 11: quint 91 kaleids"

# Update numCues for fiddlers (row 6): 34 -> 41
$ws.Range("B6").Value = 41

# Update the last-selected cell to match (cosmetic)
$ws.Range("G8").Select()
